# Edit script implementing the commit "Updated code smells and patterns
# files":
#
#   1) The "Class location" path built from several spell-check-wrapped
#      runs
#          "ganttproject" + "/" + "ganttproject" + "/" + "src" +
#          "/main/java/org/" + "imgscalr" + "/Scalr.java"
#        = "ganttproject/ganttproject/src/main/java/org/imgscalr/Scalr.java"
#      is collapsed into a single run reading
#          "/ganttproject/src/main/java/org/imgscalr/Scalr.java"
#
#   2) and 3) Two separate "Class location" paragraphs whose single run
#      reads
#          "ganttproject/biz.ganttproject.core/src/main/java/org/w3c/util/DateParser.java"
#      have the redundant leading "ganttproject" removed, leaving
#          "/biz.ganttproject.core/src/main/java/org/w3c/util/DateParser.java"

$d = $word.ActiveDocument

function Force-RunBoundary($pos) {
    # Nudge the formatting right at the given boundary (toggle bold on the
    # character immediately before it, then revert). This forces the
    # rendering engine to keep the run ending at $pos distinct from the
    # run starting at $pos instead of silently merging them together when
    # they happen to share identical run formatting.
    if ($pos -gt 0) {
        $boundary = $d.Range($pos - 1, $pos)
        $boundary.Font.Bold = 1
        $boundary.Font.Bold = 0
    }
}

# --- 1) Scalr.java: collapse the whole multi-run / proofErr-wrapped path
#        into a single run with the corrected text. The leading space
#        (belonging to the preceding run) is included in the match so
#        that the now-orphaned proofErr markers get cleaned up too; the
#        run-boundary nudge right after that space then keeps the space
#        in its own run rather than letting it merge into the new one. ---
$oldScalr = " ganttproject/ganttproject/src/main/java/org/imgscalr/Scalr.java"
$newScalr = " /ganttproject/src/main/java/org/imgscalr/Scalr.java"

$full = $d.Content.Text
$idx = $full.IndexOf($oldScalr)
if ($idx -lt 0) {
    throw "Could not find expected Scalr.java path text"
}
$rng = $d.Range($idx, $idx + $oldScalr.Length)
if ($rng.Text -ne $oldScalr) {
    throw "Unexpected text at Scalr.java match position: [$($rng.Text)]"
}
$rng.Text = $newScalr
Force-RunBoundary ($idx + 1)

# --- 2) & 3) DateParser.java: remove the redundant leading "ganttproject"
#        (12 characters) from each of the two occurrences, keeping the
#        remaining text (and its single run) otherwise untouched. ---
$dateParserNeedle = "ganttproject/biz.ganttproject.core/src/main/java/org/w3c/util/DateParser.java"

for ($i = 0; $i -lt 2; $i++) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($dateParserNeedle)
    if ($idx -lt 0) {
        throw "Could not find expected DateParser.java path text (occurrence $i)"
    }

    $rng = $d.Range($idx, $idx + 12)
    if ($rng.Text -ne "ganttproject") {
        throw "Unexpected text at DateParser.java match position: [$($rng.Text)]"
    }
    $rng.Text = ""
    Force-RunBoundary $idx
}
